# Scoreboard.xlsx update
# - Swap the "Team" (LAG N) and "Name" columns on the Score sheet (A<->B, rows 2-19)
# - Update 4 team rosters (member name changes)
# - Move the active selection to E10

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Score")

# Final team -> roster mapping (after swap + roster updates)
$teams = @(
    @{ Row = 2;  Team = "LAG 1";  Name = "Anita, Alina, Marcus, Mats" },
    @{ Row = 3;  Team = "LAG 2";  Name = "Lea, Elisabeth, Daniel, William" },
    @{ Row = 4;  Team = "LAG 3";  Name = "Marianne, Victoria, Alberto, Alekander" },
    @{ Row = 5;  Team = "LAG 4";  Name = "Elise, Kristine, Christer, Tarik" },
    @{ Row = 6;  Team = "LAG 5";  Name = "Tomine, Cecilie, Marcus, Rakan" },
    @{ Row = 7;  Team = "LAG 6";  Name = "Stine, Thea, Dario, Endre" },
    @{ Row = 8;  Team = "LAG 7";  Name = "Mari, Miriam, Bettine, Heine" },
    @{ Row = 9;  Team = "LAG 8";  Name = "Mathilde, Ida, Ivan, Snorre" },
    @{ Row = 10; Team = "LAG 9";  Name = "Regine, Sigrid, Emil, Lars" },
    @{ Row = 11; Team = "LAG 10"; Name = "Marte, Frida, Paal, Mathias" },
    @{ Row = 12; Team = "LAG 11"; Name = "Marte, Elise, Jan, Åsmund" },
    @{ Row = 13; Team = "LAG 12"; Name = "Solveig, Hedda, David, Ivan" },
    @{ Row = 14; Team = "LAG 13"; Name = "Lena, Emilie, Odd, Arne" },
    @{ Row = 15; Team = "LAG 14"; Name = "Youmna, Julie, Marius, Simone" },
    @{ Row = 16; Team = "LAG 15"; Name = "Irma, Mary, Karl, Martin" },
    @{ Row = 17; Team = "LAG 16"; Name = "Monica, Julianne, Morten, John" },
    @{ Row = 18; Team = "LAG 17"; Name = "Marianna, Sara, Carlos, Mikus" },
    @{ Row = 19; Team = "LAG 18"; Name = "Michael, Alfred, Lise, Renate" }
)

foreach ($t in $teams) {
    $ws.Cells.Item($t.Row, 1).Value = $t.Team
    $ws.Cells.Item($t.Row, 2).Value = $t.Name
}

$ws.Range("E10").Select()
